# The upstream source removed the record:
#   Caso=-560, Direccion="Pinzon 1590", OT=809098712, Proveedor Asignado="PEBCOM"
# from the field-work tracker. That record lived both in the master "General"
# sheet and in the filtered "PEBCOM" sheet (each provider has its own
# auto-generated filtered view of the same rows). Deleting the row shifts
# every following row up by one in both sheets, which is exactly what the
# diff shows (dimension A1:P394 -> A1:P393 on General, A1:P85 -> A1:P84 on
# PEBCOM, with every row below the deleted one re-numbered).

$wb = $excel.ActiveWorkbook

function Remove-CaseRow($SheetName, $CaseId, $Address) {
    $ws = $wb.Worksheets.Item($SheetName)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $caseVal = $ws.Cells.Item($r, 1).Value()
        $addrVal = $ws.Cells.Item($r, 3).Value()
        if (($caseVal -eq $CaseId) -and ($addrVal -eq $Address)) {
            $ws.Rows.Item($r).Delete()
            break
        }
    }
}

Remove-CaseRow "General" "-560" "Pinzon 1590"
Remove-CaseRow "PEBCOM" "-560" "Pinzon 1590"
